$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.919.05"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "3.890.53"
$ws.Range("E3").Value = "  +2.60%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "469.86"
$ws.Range("E5").Value = "  +9.75%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.04"
$ws.Range("E6").Value = "  +3.59%  "

$ws.Range("E7").Value = "  +0.56%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.735"
$ws.Range("E9").Value = "  +0.64%  "

$ws.Range("E10").Value = "  +8.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000336"
$ws.Range("E11").Value = "  +9.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.89"
$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").Value = "4.512.04"
$ws.Range("E13").Value = "  +3.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.32"
$ws.Range("E14").Value = "  -1.03%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.03"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").Value = "3.923.82"
$ws.Range("E16").Value = "  +4.28%  "

$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.87"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("E19").Value = "  +3.60%  "

$ws.Range("D20").Value = "67.138.38"
$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "429.53"
$ws.Range("E21").Value = "  +6.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.71"
$ws.Range("E22").Value = "  -2.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.32"
$ws.Range("E23").Value = "  +3.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.46"
$ws.Range("E24").Value = "  +4.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "38.48"
$ws.Range("E25").Value = "  +5.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.53"
$ws.Range("E26").Value = "  +7.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.76"
$ws.Range("E27").Value = "  +6.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  +1.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("E29").Value = "  -1.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "729.29"
$ws.Range("E30").Value = "  +3.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.68"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("E33").Value = "  +0.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.86"
$ws.Range("E34").Value = "  +4.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.156"
$ws.Range("E35").Value = "  +5.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.50"
$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.19%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0780"
$ws.Range("E38").Value = "  +15.39%  "

$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("E39").Value = "  -4.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0475"
$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.11"
$ws.Range("E41").Value = "  +8.49%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.141"
$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.336"
$ws.Range("E44").Value = "  +4.47%  "

$ws.Range("E45").Value = "  +5.64%  "

$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  +4.70%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  -4.81%  "

$ws.Range("E48").Value = "  +0.69%  "

$ws.Range("E49").Value = "  -1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.88"
$ws.Range("E50").Value = "  +2.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.56"
$ws.Range("E51").Value = "  +1.18%  "
